$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 291. This pushes the existing
# rows 291-302 down to 293-304 (matching the rest of the diff, which shows
# every row from 291 onward shifted down by exactly two positions with its
# original content intact).
$ws.Rows("291:292").Insert()

# Fill the two newly inserted rows with the new weekly price records.
$ws.Range("A291").Value = 11
$ws.Range("B291").Value = "Vega Monumental Concepción"
$ws.Range("C291").Value = "Bíobío"
$ws.Range("D291").Value = 44706
$ws.Range("E291").Value = 8
$ws.Range("F291").Value = 100112002
$ws.Range("G291").Value = "Pimiento"
$ws.Range("H291").Value = "Cuatro cascos verde"
$ws.Range("I291").Value = "Primera"
$ws.Range("J291").Value = 100
$ws.Range("K291").Value = 22000
$ws.Range("L291").Value = 24000
$ws.Range("M291").Value = 23000
$ws.Range("N291").Value = "`$/caja 18 kilos"
$ws.Range("O291").Value = "Provincia de Limarí"
$ws.Range("P291").Value = 1278
$ws.Range("Q291").Value = 18
$ws.Range("R291").Value = "Hortaliza"

$ws.Range("A292").Value = 11
$ws.Range("B292").Value = "Vega Monumental Concepción"
$ws.Range("C292").Value = "Bíobío"
$ws.Range("D292").Value = 44706
$ws.Range("E292").Value = 8
$ws.Range("F292").Value = 100112002
$ws.Range("G292").Value = "Pimiento"
$ws.Range("H292").Value = "Morrón rojo"
$ws.Range("I292").Value = "Primera"
$ws.Range("J292").Value = 100
$ws.Range("K292").Value = 35000
$ws.Range("L292").Value = 36000
$ws.Range("M292").Value = 35500
$ws.Range("N292").Value = "`$/caja 18 kilos"
$ws.Range("O292").Value = "Provincia de Limarí"
$ws.Range("P292").Value = 1972
$ws.Range("Q292").Value = 18
$ws.Range("R292").Value = "Hortaliza"
